$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.251.31"
$ws.Range("E2").Value = "  -0.17%  "
$ws.Range("D3").Value = "1.684.06"
$ws.Range("E3").Value = "  +0.40%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "217.15"
$ws.Range("E5").Value = "  -0.21%  "
$ws.Range("D6").Value = "0.5295"
$ws.Range("E6").Value = "  +0.75%  "
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("D8").Value = "0.2704"
$ws.Range("E8").Value = "  +0.66%  "
$ws.Range("D9").Value = "0.06390"
$ws.Range("E9").Value = "  -1.23%  "
$ws.Range("D10").Value = "21.59"
$ws.Range("E10").Value = "  -1.27%  "
$ws.Range("D11").Value = "0.07678"
$ws.Range("E11").Value = "  +2.30%  "
$ws.Range("D12").Value = "1.677.28"
$ws.Range("E12").Value = "  -1.63%  "
$ws.Range("D13").Value = "4.517"
$ws.Range("E13").Value = "  +0.09%  "
$ws.Range("D14").Value = "0.5775"
$ws.Range("E14").Value = "  +0.11%  "
$ws.Range("D15").Value = "0.000008343"
$ws.Range("E15").Value = "  -1.45%  "
$ws.Range("D16").Value = "66.55"
$ws.Range("E16").Value = "  +2.91%  "
$ws.Range("D17").Value = "26.278.98"
$ws.Range("E17").Value = "  -0.16%  "
$ws.Range("D18").Value = "1.008"
$ws.Range("E18").Value = "  -0.11%  "
$ws.Range("D19").Value = "4.883"
$ws.Range("E19").Value = "  -0.65%  "
$ws.Range("D20").Value = "10.84"
$ws.Range("E20").Value = "  -0.14%  "
$ws.Range("D21").Value = "189.79"
$ws.Range("E21").Value = "  +0.05%  "
$ws.Range("D22").Value = "6.238"
$ws.Range("E22").Value = "  +0.91%  "
$ws.Range("D23").Value = "1.009"
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("D24").Value = "149.07"
$ws.Range("E24").Value = "  +2.88%  "
$ws.Range("D25").Value = "7.807"
$ws.Range("E25").Value = "  +0.10%  "
$ws.Range("D26").Value = "0.1264"
$ws.Range("E26").Value = "  +0.32%  "
$ws.Range("D27").Value = "15.74"
$ws.Range("E27").Value = "  +0.05%  "
$ws.Range("D28").Value = "0.06267"
$ws.Range("E28").Value = "  -2.41%  "
$ws.Range("D29").Value = "1.373"
$ws.Range("E29").Value = "  +0.67%  "
$ws.Range("D30").Value = "1.322"
$ws.Range("E30").Value = "  +0.06%  "
$ws.Range("D31").Value = "3.584"
$ws.Range("E31").Value = "  +0.07%  "
$ws.Range("D32").Value = "3.565"
$ws.Range("E32").Value = "  -0.45%  "
$ws.Range("D33").Value = "1.686"
$ws.Range("E33").Value = "  +1.86%  "
$ws.Range("D34").Value = "1.025"
$ws.Range("E34").Value = "  -0.16%  "
$ws.Range("D35").Value = "0.6151"
$ws.Range("E35").Value = "  -0.41%  "
$ws.Range("D36").Value = "2.426"
$ws.Range("E36").Value = "  +0.75%  "
$ws.Range("D37").Value = "2.757"
$ws.Range("E37").Value = "  +0.54%  "
$ws.Range("D38").Value = "6.238"
$ws.Range("E38").Value = "  -0.67%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "0.01630"
$ws.Range("E39").Value = "  +0.56%  "
$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").Value = "0.8990"
$ws.Range("E40").Value = "  +3.16%  "
$ws.Range("D41").Value = "1.106.81"
$ws.Range("E41").Value = "  -1.00%  "
$ws.Range("E42").Value = "  -0.34%  "
$ws.Range("D43").Value = "100.73"
$ws.Range("E43").Value = "  +0.22%  "
$ws.Range("D44").Value = "1.836.54"
$ws.Range("E44").Value = "  +0.54%  "
$ws.Range("E45").Value = "  +3.04%  "
$ws.Range("D46").Value = "57.58"
$ws.Range("E46").Value = "  +1.18%  "
$ws.Range("E47").Value = "  -0.13%  "
$ws.Range("D48").Value = "8.091"
$ws.Range("E48").Value = "  -0.96%  "
$ws.Range("D49").Value = "0.05280"
$ws.Range("E49").Value = "  +0.30%  "
$ws.Range("D50").Value = "0.4292"
$ws.Range("E50").Value = "  -0.16%  "
$ws.Range("D51").Value = "6.035"
$ws.Range("E51").Value = "  -0.25%  "
